$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# New homework data table (rows 4-15) on sheet "2": for a few people, list
# their birthplace / registered address / work address, joined by name.
# Filled column by column (A, then B, then C) to match the original
# authoring order of the shared-strings table.
# ---------------------------------------------------------------------------

$names  = @("Иванов И. И.", "Иванов И. И.", "Иванов И. И.", "Иванов И. И.",
            "Петров П. П.", "Васильев В. В.",
            "null", "null", "null", "null", "null", "null")

$places = @("Можга", "Казань", "Москва", "Санкт-Петербург",
            "Москва", "Белгород",
            "Уфа", "Сочи", "Киров", "Владивосток", "Рязань", "Хабаровск")

$kinds  = @("Место рождения", "По прописке", "Рабочий", "По прописке",
            "По прописке", "По прописке",
            "По прописке", "По прописке", "Рабочий", "Место рождения",
            "Рабочий", "Место рождения")

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws2.Cells.Item($i + 4, 1).Value = $names[$i]
}
for ($i = 0; $i -lt $places.Count; $i++) {
    $ws2.Cells.Item($i + 4, 2).Value = $places[$i]
}
for ($i = 0; $i -lt $kinds.Count; $i++) {
    $ws2.Cells.Item($i + 4, 3).Value = $kinds[$i]
}

# ---------------------------------------------------------------------------
# Formatting: reuse the styles already present on sheet "1"
#   - A4 style (bordered cell) for columns B/C and most of column A
#   - C13 style (left/right border only) for the repeated "null" rows (10-15)
# ---------------------------------------------------------------------------

$ws1.Range("A4").Copy()
$ws2.Range("A4:C9").PasteSpecial(-4122)
$ws2.Range("B10:C15").PasteSpecial(-4122)

$ws1.Range("C13").Copy()
$ws2.Range("A10:A15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Rows with two-line wrapped text need the taller row height.
$ws2.Rows.Item(4).RowHeight = 30
$ws2.Rows.Item(7).RowHeight = 30
$ws2.Rows.Item(13).RowHeight = 30
$ws2.Rows.Item(15).RowHeight = 30

# ---------------------------------------------------------------------------
# Restore the active cell selection on sheet "2" to match the saved state.
# ---------------------------------------------------------------------------
$ws2.Range("E22").Select()
